$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the schedule grid: move the scheduled shift block from rows 6-7
# (Monday-Friday, columns B-F) down to row 30 (Monday-Friday, columns B-F).
$ws.Range("B6:F7").Value = 0
$ws.Range("B30:F30").Value = 1

# Update the view state: scroll the window so row 10 is the top-visible row,
# then select the newly-scheduled range so it becomes the active selection.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B30:F30").Select()
